{"js": "// Helper: find the paragraph containing `anchorText` and replace the\n// paragraph's whole text with `newText`. Using paragraph.getRange(\"Whole\")\n// (rather than the search-hit range itself) lets the run(s) collapse\n// cleanly into a single run, mirroring what Word does when the runs in a\n// paragraph are merged/retyped.\nasync function replaceParagraphText(anchorText, newText) {\n  const results = context.document.body.search(anchorText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  const paragraphs = results.items[0].paragraphs;\n  paragraphs.load(\"items\");\n  await context.sync();\n\n  const wholeParagraphRange = paragraphs.items[0].getRange(\"Whole\");\n  wholeParagraphRange.insertText(newText, \"Replace\");\n  await context.sync();\n}\n\n// 1) \"How to ... configure your setup for development of ... the\n//    CDInterface ... Module\" -> runs merged into one (text unchanged).\nawait replaceParagraphText(\n  \"How to configure your setup for development of the CDInterfaceModule\",\n  \"How to configure your setup for development of the CDInterfaceModule\"\n);\n\n// 2) \"Please refer to the ... following documents for related\n//    information\" -> runs merged into one (text unchanged).\nawait replaceParagraphText(\n  \"Please refer to the following documents for related information\",\n  \"Please refer to the following documents for related information\"\n);\n\n// 3) \"CDInterfaceModule ... User ... Guide for more details on how to\n//    ... use ... the application.\" -> runs merged into one (text\n//    unchanged).\nawait replaceParagraphText(\n  \"CDInterfaceModule User Guide for more details on how to use the application.\",\n  \"CDInterfaceModule User Guide for more details on how to use the application.\"\n);\n\n// 4) \"CDInterfaceModule ... Design ... for more details on how to ...\n//    the code works within t ... he application.\" -> runs merged into\n//    one (text unchanged).\nawait replaceParagraphText(\n  \"CDInterfaceModule Design for more details on how to the code works within the application.\",\n  \"CDInterfaceModule Design for more details on how to the code works within the application.\"\n);\n\n// 5) deploy.ps1 paragraph: actual wording change (added \"all\"/\"cleanall\"\n//    and an extra trailing sentence).\nawait replaceParagraphText(\n  \"The deploy.ps1 script has an option which will remove the installed version of the module, -cleanonly.\",\n  \"The deploy.ps1 script has an option which will remove all the installed versions of the module, -cleanallonly. The -cleanonly option will just remove the current version.\"\n);\n\n// 6) Mark each inline picture's range as \"do not spell check\" (adds\n//    <w:noProof/> to the run properties of the run hosting the drawing),\n//    matching the three <w:noProof/> additions in the diff.\nconst inlinePictures = context.document.body.inlinePictures;\ninlinePictures.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < inlinePictures.items.length; i++) {\n  inlinePictures.items[i].getRange().hasNoProofing = true;\n}\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Replace the text of the first range matching $findText with $replaceText.\n# Using Find/Execute's own Replace (instead of assigning Range.Text\n# directly) guarantees the matched run(s) are rewritten even when the\n# replacement text is identical to the concatenation of the old runs -\n# which is exactly what happens for the paragraphs below whose multiple\n# runs are simply being merged back into a single run.\nfunction Replace-ParagraphText($findText, $replaceText) {\n    $find = $d.Content.Find\n    $find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 1) | Out-Null\n}\n\n# 1) \"How to \" + \"configure your setup for development of \" +\n#    \"the CDInterface\" + \"Module\" -> merged into a single run (text\n#    itself is unchanged).\nReplace-ParagraphText \"How to configure your setup for development of the CDInterfaceModule\" \"How to configure your setup for development of the CDInterfaceModule\"\n\n# 2) \"Please refer to the \" + \"following documents for related\n#    information\" -> merged into a single run (text unchanged).\nReplace-ParagraphText \"Please refer to the following documents for related information\" \"Please refer to the following documents for related information\"\n\n# 3) \"CDInterfaceModule \" + \"User\" + \" Guide for more details on how to \"\n#    + \"use\" + \" the application.\" -> merged into a single run (text\n#    unchanged).\nReplace-ParagraphText \"CDInterfaceModule User Guide for more details on how to use the application.\" \"CDInterfaceModule User Guide for more details on how to use the application.\"\n\n# 4) \"CDInterfaceModule \" + \"Design\" + \" for more details on how to \" +\n#    \"the code works within t\" + \"he application.\" -> merged into a\n#    single run (text unchanged).\nReplace-ParagraphText \"CDInterfaceModule Design for more details on how to the code works within the application.\" \"CDInterfaceModule Design for more details on how to the code works within the application.\"\n\n# 5) deploy.ps1 paragraph: actual wording change (added \"all\"/\"cleanall\"\n#    and an extra trailing sentence).\nReplace-ParagraphText \"The deploy.ps1 script has an option which will remove the installed version of the module, -cleanonly.\" \"The deploy.ps1 script has an option which will remove all the installed versions of the module, -cleanallonly. The -cleanonly option will just remove the current version.\"\n\n# 6) Mark each inline picture as \"do not spell check\" (adds\n#    <w:noProof/> to the run properties of the run hosting the drawing),\n#    matching the three <w:noProof/> additions in the diff.\n$shapes = $d.InlineShapes\nfor ($i = 1; $i -le $shapes.Count; $i++) {\n    $shapes.Item($i).Range.NoProofing = $true\n}\n"}
